# Petty cash book update - 15 Feb 2021, end of day update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New opening balance carried forward
$ws.Range("E2").Value = 854500

# New transaction entry for 15 Feb 2021 (Wages Expense category kept, no amount yet)
$ws.Range("A3").Value = 44242
$ws.Range("D3").Clear()

# Clear out all the old transaction detail rows (date/description/debit/credit),
# leaving only the running-balance formulas in column E
$ws.Range("A4:D33").Clear()

# Update the view: scroll position and active selection
$ws.Range("B5").Select() | Out-Null
